$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; existing rows 14-48 shift down to 15-49.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly price record.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 45037
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100114007
$ws.Range("G14").Value = "Jengibre"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 16000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16400
$ws.Range("N14").Value = "`$/caja 13 kilos"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1262
$ws.Range("Q14").Value = 13
$ws.Range("R14").Value = "Hortaliza"
